# Apply cryptos list update per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '44.204.65'
$ws.Cells.Item(2, 5).Value = '  +1.90%  '

# Row 3
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.357.81'
$ws.Cells.Item(3, 5).Value = '  -0.03%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.680'
$ws.Cells.Item(5, 5).Value = '  +4.61%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '238.84'
$ws.Cells.Item(6, 5).Value = '  +2.50%  '

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '75.18'
$ws.Cells.Item(7, 5).Value = '  +10.82%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.07%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +20.46%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.101'
$ws.Cells.Item(10, 5).Value = '  +5.49%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '30.76'
$ws.Cells.Item(11, 5).Value = '  +16.71%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.21%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.707.65'
$ws.Cells.Item(13, 5).Value = '  +0.16%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '16.98'
$ws.Cells.Item(14, 5).Value = '  +7.84%  '

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.91'
$ws.Cells.Item(15, 5).Value = '  +10.08%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.910'
$ws.Cells.Item(16, 5).Value = '  +7.86%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.367.72'
$ws.Cells.Item(17, 5).Value = '  +0.27%  '

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '44.276.55'
$ws.Cells.Item(18, 5).Value = '  +2.12%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +3.99%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '78.27'
$ws.Cells.Item(20, 5).Value = '  +5.70%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.52'
$ws.Cells.Item(21, 5).Value = '  +4.07%  '

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '255.73'
$ws.Cells.Item(22, 5).Value = '  +2.50%  '

# Row 23
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$ws.Cells.Item(23, 5).Value = '  -0.10%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'WEMIXToken'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.79'
$ws.Cells.Item(24, 5).Value = '  -4.90%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +2.62%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Toncoin'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.39'
$ws.Cells.Item(26, 5).Value = '  +4.86%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.43'
$ws.Cells.Item(27, 5).Value = '  +5.22%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.92%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '174.14'
$ws.Cells.Item(29, 5).Value = '  +1.10%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +3.30%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.42%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +5.19%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +7.25%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +3.73%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +3.48%  '

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.88'
$ws.Cells.Item(36, 5).Value = '  +6.61%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -2.43%  '

# Row 38
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.47'
$ws.Cells.Item(38, 5).Value = '  -0.82%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +6.08%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.22'
$ws.Cells.Item(40, 5).Value = '  +5.06%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.06%  '

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.88'
$ws.Cells.Item(42, 5).Value = '  -1.19%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +3.78%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  +4.17%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Algorand'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.187'
$ws.Cells.Item(45, 5).Value = '  +12.10%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'ARBITRUM'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.17'
$ws.Cells.Item(46, 5).Value = '  +0.14%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '99.52'
$ws.Cells.Item(47, 5).Value = '  +0.64%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +0.24%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +5.14%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.444.96'
$ws.Cells.Item(50, 5).Value = '  -0.32%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000204'
$ws.Cells.Item(51, 5).Value = '  +1.62%  '
